# Se adiciona el control de fecha de inicio y final de autorizacion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new columns
$ws.Range("J1").Value = "StartTime"
$ws.Range("K1").Value = "EndTime"

# Column widths matching the new columns (ColumnWidth setter adds a fixed
# 5/6-character padding when the value round-trips through the stored
# <col width=.../> attribute, so compensate to land on the exact target).
$ws.Columns.Item(10).ColumnWidth = 19.5 - (5/6)
$ws.Columns.Item(11).ColumnWidth = 18 - (5/6)

# Date/time values (stored as Excel serial date-times)
$ws.Range("J2").Value = 44018.999988425923
$ws.Range("K2").Value = 44171.999988425923
$ws.Range("J3").Value = 44018.999988425923
$ws.Range("K3").Value = 44171.999988425923

# Apply font + number format to the new date cells in one pass (matches
# target style: a single extra cell style combining font size and format)
$dateRng = $ws.Range("J2:K3")
$dateRng.Font.Size = 10.5
$dateRng.NumberFormat = "yyyy/mm/dd\ h:mm:ss"

# Update selection to match target workbook state
$ws.Range("K4").Select()

$wb.Save()
